# Auto-generated edit script: updates market-price-derived columns (H:N)
# on the Leviathan_Profits workbook sheets, per the authoritative diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 570.4545000000001
$ws.Range("I11").Value = 570.4545000000001
$ws.Range("K11").Value = 570.4545000000001
$ws.Range("M11").Value = -430.4545000000001
# Row 62
$ws.Range("H62").Value = 4377.3335
$ws.Range("I62").Value = 3983.0625
$ws.Range("K62").Value = 3983.0625
$ws.Range("M62").Value = -3359.0625
# Row 65
$ws.Range("H65").Value = 4377.3335
$ws.Range("I65").Value = 3983.0625
$ws.Range("K65").Value = 19915.3125
$ws.Range("M65").Value = -16795.3125
# Row 98
$ws.Range("H98").Value = 1820.8334
$ws.Range("I98").Value = 1181.8182
$ws.Range("K98").Value = 1181.8182
$ws.Range("M98").Value = 316.1818000000001
# Row 122
$ws.Range("H122").Value = 1820.8334
$ws.Range("I122").Value = 1181.8182
$ws.Range("K122").Value = 3545.4546
$ws.Range("M122").Value = -1095.4546
# Row 137
$ws.Range("H137").Value = 1968.1904
$ws.Range("I137").Value = 1859.421
$ws.Range("K137").Value = 5578.263
$ws.Range("M137").Value = -3028.263

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1657.4762
$ws.Range("I2").Value = 1699.3334
$ws.Range("J2").Value = 1552.8334
$ws.Range("K2").Value = 1699.3334
$ws.Range("L2").Value = 1552.8334
$ws.Range("M2").Value = -1586.3334
$ws.Range("N2").Value = -1778.8334
# Row 3
$ws.Range("H3").Value = 750
$ws.Range("I3").Value = 750
$ws.Range("K3").Value = 750
$ws.Range("M3").Value = -635
# Row 32
$ws.Range("H32").Value = 40795.9
$ws.Range("I32").Value = 7551.185
$ws.Range("K32").Value = 7551.185
$ws.Range("M32").Value = -7264.185
# Row 45
$ws.Range("H45").Value = 6423.885
$ws.Range("I45").Value = 10445.333
$ws.Range("K45").Value = 10445.333
$ws.Range("M45").Value = -10068.333
# Row 57
$ws.Range("H57").Value = 5160
$ws.Range("I57").Value = 5160
$ws.Range("K57").Value = 5160
$ws.Range("M57").Value = -4676
# Row 74
$ws.Range("H74").Value = 2522.5264
$ws.Range("I74").Value = 2362.1936
$ws.Range("J74").Value = 3232.5715
$ws.Range("K74").Value = 2362.1936
$ws.Range("L74").Value = 3232.5715
$ws.Range("M74").Value = -1488.1936
$ws.Range("N74").Value = -4980.5715
# Row 77
$ws.Range("H77").Value = 2522.5264
$ws.Range("I77").Value = 2362.1936
$ws.Range("J77").Value = 3232.5715
$ws.Range("K77").Value = 11810.968
$ws.Range("L77").Value = 16162.8575
$ws.Range("M77").Value = -7442.968000000001
$ws.Range("N77").Value = -24898.8575
# Row 116
$ws.Range("H116").Value = 1657.4762
$ws.Range("I116").Value = 1699.3334
$ws.Range("J116").Value = 1552.8334
$ws.Range("K116").Value = 1699.3334
$ws.Range("L116").Value = 1552.8334
$ws.Range("M116").Value = 594.6666
$ws.Range("N116").Value = -6140.8334
# Row 122
$ws.Range("H122").Value = 3319.3
$ws.Range("I122").Value = 3347.875
$ws.Range("K122").Value = 10043.625
$ws.Range("M122").Value = -7593.625
# Row 132
$ws.Range("H132").Value = 3822.1904
$ws.Range("I132").Value = 1559.6923
$ws.Range("J132").Value = 7498.75
$ws.Range("K132").Value = 4679.0769
$ws.Range("L132").Value = 22496.25
$ws.Range("M132").Value = -2149.0769
$ws.Range("N132").Value = -27556.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1657.4762
$ws.Range("I3").Value = 1699.3334
$ws.Range("J3").Value = 1552.8334
$ws.Range("K3").Value = 1699.3334
$ws.Range("L3").Value = 1552.8334
$ws.Range("M3").Value = -1585.3334
$ws.Range("N3").Value = -1780.8334
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
# Row 134
$ws.Range("H134").Value = 831.8276
$ws.Range("I134").Value = 652.6818
$ws.Range("J134").Value = 1394.8572
$ws.Range("K134").Value = 1958.0454
$ws.Range("L134").Value = 4184.571599999999
$ws.Range("M134").Value = 576.9546
$ws.Range("N134").Value = -9254.571599999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 45411.29
$ws.Range("I31").Value = 40911.883
$ws.Range("J31").Value = 68808.2
$ws.Range("K31").Value = 40911.883
$ws.Range("L31").Value = 68808.2
$ws.Range("M31").Value = -40616.883
$ws.Range("N31").Value = -69398.2
# Row 34
$ws.Range("H34").Value = 45411.29
$ws.Range("I34").Value = 40911.883
$ws.Range("J34").Value = 68808.2
$ws.Range("K34").Value = 40911.883
$ws.Range("L34").Value = 68808.2
$ws.Range("M34").Value = -40709.883
$ws.Range("N34").Value = -69212.2
# Row 99
$ws.Range("H99").Value = 35473.5
$ws.Range("I99").Value = 46333
$ws.Range("J99").Value = 2895
$ws.Range("K99").Value = 46333
$ws.Range("L99").Value = 2895
$ws.Range("M99").Value = -44835
$ws.Range("N99").Value = -5891
# Row 105
$ws.Range("H105").Value = 2500
$ws.Range("J105").Value = 2500
$ws.Range("L105").Value = 2500
$ws.Range("N105").Value = -5994
# Row 126
$ws.Range("H126").Value = 35473.5
$ws.Range("I126").Value = 46333
$ws.Range("J126").Value = 2895
$ws.Range("K126").Value = 138999
$ws.Range("L126").Value = 8685
$ws.Range("M126").Value = -136529
$ws.Range("N126").Value = -13625
# Row 132
$ws.Range("H132").Value = 3172.4866
$ws.Range("I132").Value = 3112.625
$ws.Range("K132").Value = 9337.875
$ws.Range("M132").Value = -6807.875
# Row 134
$ws.Range("H134").Value = 2592.9443
$ws.Range("I134").Value = 2568.1765
$ws.Range("K134").Value = 7704.529500000001
$ws.Range("M134").Value = -5169.529500000001
# Row 141
$ws.Range("H141").Value = 365143.25
$ws.Range("J141").Value = 365143.25
$ws.Range("L141").Value = 365143.25
$ws.Range("N141").Value = -375503.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 827.8
$ws.Range("I5").Value = 786.125
$ws.Range("J5").Value = 994.5
$ws.Range("K5").Value = 2358.375
$ws.Range("L5").Value = 2983.5
$ws.Range("M5").Value = -2246.375
$ws.Range("N5").Value = -3207.5
# Row 135
$ws.Range("H135").Value = 827.8
$ws.Range("I135").Value = 786.125
$ws.Range("J135").Value = 994.5
$ws.Range("K135").Value = 7075.125
$ws.Range("L135").Value = 8950.5
$ws.Range("M135").Value = -4540.125
$ws.Range("N135").Value = -14020.5
# Row 139
$ws.Range("H139").Value = 5202.2
$ws.Range("I139").Value = 3009.6667
$ws.Range("J139").Value = 6141.857
$ws.Range("K139").Value = 9029.000100000001
$ws.Range("L139").Value = 18425.571
$ws.Range("M139").Value = -3889.000100000001
$ws.Range("N139").Value = -28705.571

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
# Row 70
$ws.Range("H70").Value = 5672.7
$ws.Range("J70").Value = 4871.6665
$ws.Range("L70").Value = 4871.6665
$ws.Range("N70").Value = -5411.6665
# Row 73
$ws.Range("H73").Value = 5672.7
$ws.Range("J73").Value = 4871.6665
$ws.Range("L73").Value = 4871.6665
$ws.Range("N73").Value = -6743.6665
# Row 102
$ws.Range("H102").Value = 1591.48
$ws.Range("I102").Value = 1373.3478
$ws.Range("K102").Value = 1373.3478
$ws.Range("M102").Value = 248.6522

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 4159
$ws.Range("I122").Value = 3312.9092
$ws.Range("J122").Value = 5089.7
$ws.Range("K122").Value = 9938.7276
$ws.Range("L122").Value = 15269.1
$ws.Range("M122").Value = -7488.7276
$ws.Range("N122").Value = -20169.1
# Row 132
$ws.Range("H132").Value = 6665
$ws.Range("J132").Value = 6665
$ws.Range("L132").Value = 19995
$ws.Range("N132").Value = -25055

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 97
$ws.Range("H97").Value = 36179
$ws.Range("J97").Value = 36179
$ws.Range("L97").Value = 36179
$ws.Range("N97").Value = -38161
# Row 107
$ws.Range("H107").Value = 20834034
$ws.Range("I107").Value = 691.7
$ws.Range("K107").Value = 2075.1
$ws.Range("M107").Value = -155.1000000000004
# Row 113
$ws.Range("H113").Value = 439.7143
$ws.Range("I113").Value = 215.6
$ws.Range("K113").Value = 646.8
$ws.Range("M113").Value = 1523.2
# Row 126
$ws.Range("H126").Value = 878
$ws.Range("I126").Value = 776.4783
$ws.Range("K126").Value = 2329.4349
$ws.Range("M126").Value = 140.5650999999998
# Row 132
$ws.Range("H132").Value = 4907.2354
$ws.Range("I132").Value = 5194.8667
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 15584.6001
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = -13054.6001
$ws.Range("N132").Value = -13310
# Row 136
$ws.Range("H136").Value = 1228.9231
$ws.Range("I136").Value = 736.5714
$ws.Range("K136").Value = 2209.7142
$ws.Range("M136").Value = 340.2857999999997
